# Update cryptos list (prices / volume% changes + Mantle/Kaspa row swap)
# Commit: "Updated cryptos list on Thu May 30 09:19:54 UTC 2024 with GitHub Actions"
#
# NOTE: the Price column (D) stores values as TEXT (inline strings) even
# though most look numeric (e.g. "591.48", "3.731.61"). Assigning such a
# string straight to Range.Value lets Excel auto-coerce it into a real
# number (changing both the stored type and introducing float rounding,
# e.g. 591.48 -> 591.48000000000002). To keep these as text -- matching
# the source file -- force NumberFormat="@" before writing the value,
# then reset Style back to "Normal" so no stray style index lingers on
# the cell. Columns B/C/E never look numeric to Excel (URLs, names, and
# padded "  +0.00%  " strings), so they can be assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '67.686.81'
$c.Style = 'Normal'
$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '3.734.38'
$c.Style = 'Normal'
$ws.Range("E3").Value = '  -1.92%  '
$ws.Range("E4").Value = '  -0.14%  '
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '591.48'
$c.Style = 'Normal'
$ws.Range("E5").Value = '  -1.25%  '
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '165.51'
$c.Style = 'Normal'
$ws.Range("E6").Value = '  -2.08%  '
$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '3.732.51'
$c.Style = 'Normal'
$ws.Range("E7").Value = '  -1.92%  '
$ws.Range("E8").Value = '  -0.03%  '
$c = $ws.Range("D9")
$c.NumberFormat = '@'
$c.Value = '0.518'
$c.Style = 'Normal'
$ws.Range("E9").Value = '  -2.20%  '
$ws.Range("E10").Value = '  -3.89%  '
$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '6.46'
$c.Style = 'Normal'
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("E12").Value = '  -2.40%  '
$ws.Range("E13").Value = '  -5.34%  '
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '36.01'
$c.Style = 'Normal'
$ws.Range("E14").Value = '  -2.50%  '
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '4.361.66'
$c.Style = 'Normal'
$ws.Range("E15").Value = '  -1.88%  '
$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '3.725.89'
$c.Style = 'Normal'
$ws.Range("E16").Value = '  -1.68%  '
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '67.667.45'
$c.Style = 'Normal'
$ws.Range("E17").Value = '  -0.53%  '
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '18.28'
$c.Style = 'Normal'
$ws.Range("E18").Value = '  -0.69%  '
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '7.03'
$c.Style = 'Normal'
$ws.Range("E19").Value = '  -5.76%  '
$ws.Range("E20").Value = '  -0.50%  '
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '10.65'
$c.Style = 'Normal'
$ws.Range("E21").Value = '  -2.29%  '
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '467.02'
$c.Style = 'Normal'
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '0.699'
$c.Style = 'Normal'
$ws.Range("E23").Value = '  -4.31%  '
$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '82.71'
$c.Style = 'Normal'
$ws.Range("E24").Value = '  -0.86%  '
$ws.Range("E25").Value = '  -10.88%  '
$ws.Range("E26").Value = '  -5.84%  '
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '11.97'
$c.Style = 'Normal'
$ws.Range("E27").Value = '  -2.20%  '
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '10.08'
$c.Style = 'Normal'
$ws.Range("E28").Value = '  -2.23%  '
$ws.Range("E29").Value = '  +0.02%  '
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '3.879.61'
$c.Style = 'Normal'
$ws.Range("E30").Value = '  -1.89%  '
$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '2.77'
$c.Style = 'Normal'
$ws.Range("E31").Value = '  -5.20%  '
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '7.36'
$c.Style = 'Normal'
$ws.Range("E32").Value = '  -4.93%  '
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '2.22'
$c.Style = 'Normal'
$ws.Range("E33").Value = '  -3.10%  '
$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '29.58'
$c.Style = 'Normal'
$ws.Range("E34").Value = '  -4.06%  '
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '9.01'
$c.Style = 'Normal'
$ws.Range("E35").Value = '  -3.46%  '
$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '3.686.29'
$c.Style = 'Normal'
$ws.Range("E36").Value = '  -2.30%  '
$ws.Range("E37").Value = '  -5.66%  '
$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '3.42'
$c.Style = 'Normal'
$ws.Range("E38").Value = '  -10.33%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '0.138'
$c.Style = 'Normal'
$ws.Range("E39").Value = '  -1.46%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '0.991'
$c.Style = 'Normal'
$ws.Range("E40").Value = '  -2.30%  '
$ws.Range("E41").Value = '  -3.95%  '
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("E43").Value = '  -0.06%  '
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '0.305'
$c.Style = 'Normal'
$ws.Range("E44").Value = '  -4.34%  '
$ws.Range("E45").Value = '  -3.23%  '
$ws.Range("E46").Value = '  -3.49%  '
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '45.28'
$c.Style = 'Normal'
$ws.Range("E47").Value = '  -2.66%  '
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '390.81'
$c.Style = 'Normal'
$ws.Range("E48").Value = '  -4.19%  '
$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '143.33'
$c.Style = 'Normal'
$ws.Range("E49").Value = '  +0.94%  '
$ws.Range("E50").Value = '  -3.69%  '
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '25.05'
$c.Style = 'Normal'
$ws.Range("E51").Value = '  -1.68%  '